# Link Between MIFCRIANCA And MIFCRIANCA_VISIT

$wb = $excel.ActiveWorkbook

# NOTE: new shared strings are interned in the order cells are written, and the
# saved workbook expects them in the order MIFCRIANCA, hash-formula, CHILD - so
# the writes below are deliberately interleaved across the two sheets to match.

$wsChoices = $wb.Worksheets.Item("choices")
$wsSurvey = $wb.Worksheets.Item("survey")

# --- choices sheet: add a new "forms" choice row for MIFCRIANCA / CHILD ---
$wsChoices.Range("A11").Value = "forms"
$wsChoices.Range("B11").Value = "MIFCRIANCA"

# --- survey sheet: add the MIFCRIANCA branch entries (mirrors the MIFCRIANCA_VISIT block above it) ---
$wsSurvey.Range("A35").Value = "MIFCRIANCA"
# Leading "'" is Excel's quote-prefix marker and is swallowed on assignment, so
# double it to keep a single literal leading apostrophe in the stored text.
$wsSurvey.Range("B36").Value = "''?' + odkSurvey.getHashString('MIFCRIANCA')"
$wsSurvey.Range("E36").Value = "external_link"
$wsSurvey.Range("G36").Value = "Open form"
$wsSurvey.Range("C37").Value = "exit section"

# --- back to choices sheet for the CHILD data_value / display.title.text columns ---
$wsChoices.Range("C11").Value = "CHILD"
$wsChoices.Range("D11").Value = "CHILD"

# --- restore view/selection state to match the saved workbook ---
$wsSurvey.Activate()
$wsSurvey.Range("D27").Select()

$wsChoices.Activate()
$wsChoices.Range("C16").Select()
